$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I20").Value = 0.2049479017366116
$ws.Range("J20").Value = 0.263184587419064
$ws.Range("K20").Value = 0.2449524956484287
$ws.Range("L20").Value = 2.431376744545068
